$wb = $excel.ActiveWorkbook

# Update the "zh-cn" report sheet: refresh the handoff/handback datetimes
# for the db01e9df-a431-4992-9bca-71a2ae7f91f0 row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-18 02:42:19"
$wsZhCn.Range("G4").Value = "2016-01-18 02:43:22"

# Update the "de-de" report sheet: refresh the handoff/handback datetimes
# for the db01e9df-a431-4992-9bca-71a2ae7f91f0 row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-18 02:42:33"
$wsDeDe.Range("G4").Value = "2016-01-18 02:43:46"
